$d = $word.ActiveDocument

# The date line reads "Күні: 19.02.2024ж" and must become "Күні: 20.02.2024ж".
# The two runs holding "1" and "9" are merged into a single "20" run, so do
# a plain text replace of "19.02.2024" -> "20.02.2024" (unique in the document).
$d.Content.Find.Execute("19.02.2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "20.02.2024", 2)
